# Update "paises" / provincias data: refresh country case counters and the
# "last updated" timestamp, then re-sort the table by total cases (column B)
# descending - same as the source sheet always is - so any country whose
# updated total now overtakes its neighbour changes rows naturally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last updated" banner in A1 -----------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 19:39"

# --- 2. Update the per-country counters ------------------------------------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Estados Unidos
$ws.Range("B4").Value = 4341491
$ws.Range("C4").Value = 25782
$ws.Range("D4").Value = 2072518
$ws.Range("E4").Value = 2119372
$ws.Range("G4").Value = 203
$ws.Range("H4").Value = 149601

# India
$ws.Range("B6").Value = 1435213
$ws.Range("C6").Value = 49719
$ws.Range("D6").Value = 917234
$ws.Range("E6").Value = 485170
$ws.Range("G6").Value = 713
$ws.Range("H6").Value = 32809

# Chile
$ws.Range("B11").Value = 345790
$ws.Range("C11").Value = 2198
$ws.Range("D11").Value = 318095
$ws.Range("E11").Value = 18583
$ws.Range("G11").Value = 92
$ws.Range("H11").Value = 9112

# Turquia
$ws.Range("B19").Value = 226100
$ws.Range("C19").Value = 927
$ws.Range("D19").Value = 209487
$ws.Range("E19").Value = 11000
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 5613

# Israel
$ws.Range("B41").Value = 61764
$ws.Range("C41").Value = 1086
$ws.Range("D41").Value = 27014
$ws.Range("E41").Value = 34282
$ws.Range("G41").Value = 11
$ws.Range("H41").Value = 468

# Argelia
$ws.Range("B60").Value = 27357
$ws.Range("C60").Value = 593
$ws.Range("D60").Value = 18088
$ws.Range("E60").Value = 8114
$ws.Range("G60").Value = 9
$ws.Range("H60").Value = 1155

# Irlanda
$ws.Range("B61").Value = 25881
$ws.Range("C61").Value = 12
$ws.Range("E61").Value = 753

# Marruecos (currently row 66; will move above Uzbekistan once sorted because
# its new total, 20278, now exceeds Uzbekistan's 20226)
$ws.Range("B66").Value = 20278
$ws.Range("C66").Value = 633
$ws.Range("D66").Value = 16438
$ws.Range("E66").Value = 3527
$ws.Range("G66").Value = 8
$ws.Range("H66").Value = 313

# Libano
$ws.Range("B105").Value = 3750
$ws.Range("C105").Value = 168
$ws.Range("E105").Value = 2028
$ws.Range("G105").Value = 4
$ws.Range("H105").Value = 51

# Tunez
$ws.Range("B137").Value = 1452
$ws.Range("C137").Value = 9
$ws.Range("D137").Value = 1142

# Uganda
$ws.Range("B145").Value = 1115
$ws.Range("C145").Value = 12
$ws.Range("E145").Value = 131

# --- 3. Re-sort the country table by total cases (desc), as the sheet is
#        always kept, so rows whose ranking changed settle into place -------
$dataRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$dataRange.Sort($keyRange, 2)
